$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure mobile number / college code columns stay text (preserve leading zeros)
$ws.Range("D2:D4").NumberFormat = "@"
$ws.Range("G2:G4").NumberFormat = "@"

# Row 2: replace Devam Narkar's data with SHAIKH MG's data
$ws.Range("B2").Value = "SHAIKH MG"
$ws.Range("C2").Value = "dec612owner@gtu.edu.in"
$ws.Range("D2").Value = "9173971588"
$ws.Range("E2").Value = "DR. S. & S. S. GHANDHY COLLEGE OF ENGINEERING & TECHNOLOGY"
$ws.Range("F2").Value = "surat"
$ws.Range("G2").Value = "612"
$ws.Range("H2").Value = "Gujarat"

# Row 3: Kalpesh Shah
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Kalpesh Shah"
$ws.Range("C3").Value = "meet_kalpesh@yahoo.co.in"
$ws.Range("D3").Value = "9427062268"
$ws.Range("E3").Value = "A. D. PATEL INSTITUTE OF TECHNOLOGY"
$ws.Range("F3").Value = "anand"
$ws.Range("G3").Value = "001"
$ws.Range("H3").Value = "Gujarat"

# Row 4: VAIDYA CHIRAYU
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "VAIDYA CHIRAYU"
$ws.Range("C4").Value = "cpv.fetr@gmail.com"
$ws.Range("D4").Value = "9737914668"
$ws.Range("E4").Value = "FACULTY OF ENGINEERING ,TECHNOLOGY AND RESEARCH. BARDOLI"
$ws.Range("F4").Value = "surat"
$ws.Range("G4").Value = "084"
$ws.Range("H4").Value = "Gujarat"
